# Apply "Update latest output (run 182)" changes
$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Schedule" ---
$ws1 = $wb.Worksheets.Item("Schedule")

# Remove the old row 3 (schedule block that no longer exists after re-optimisation)
$ws1.Rows.Item(3).Delete()

# Update row 2 with the new schedule values
$ws1.Range("A2").Value2 = 46045
$ws1.Range("B2").Value2 = 46045.66666666666
$ws1.Range("C2").Value2 = 16
$ws1.Range("D2").Value2 = 60.48
$ws1.Range("E2").Value2 = 1325.030382
$ws1.Range("F2").Value2 = 21.90857113095239

# --- Sheet 2: "Detailed" ---
$ws2 = $wb.Worksheets.Item("Detailed")

# Grow the table by one row (48 data rows -> 49 data rows), appending after the
# existing last row so the new row inherits the same number formats/styles.
$ws2.Rows.Item(49).Insert()

$arr = New-Object 'object[,]' 48,5
$arr[0,0] = 46045
$arr[0,1] = 57.06
$arr[0,2] = "historical"
$arr[0,3] = 46045
$arr[0,4] = "ON"
$arr[1,0] = 46045.02083333334
$arr[1,1] = 57.06
$arr[1,2] = "historical"
$arr[1,3] = 46045
$arr[1,4] = "ON"
$arr[2,0] = 46045.04166666666
$arr[2,1] = 57.06
$arr[2,2] = "historical"
$arr[2,3] = 46045
$arr[2,4] = "ON"
$arr[3,0] = 46045.0625
$arr[3,1] = 56.98
$arr[3,2] = "historical"
$arr[3,3] = 46045
$arr[3,4] = "ON"
$arr[4,0] = 46045.08333333334
$arr[4,1] = 56.97994
$arr[4,2] = "historical"
$arr[4,3] = 46045
$arr[4,4] = "ON"
$arr[5,0] = 46045.10416666666
$arr[5,1] = 56.97994
$arr[5,2] = "forecast"
$arr[5,3] = 46045
$arr[5,4] = "ON"
$arr[6,0] = 46045.125
$arr[6,1] = 57.06
$arr[6,2] = "forecast"
$arr[6,3] = 46045
$arr[6,4] = "ON"
$arr[7,0] = 46045.14583333334
$arr[7,1] = 57.06
$arr[7,2] = "forecast"
$arr[7,3] = 46045
$arr[7,4] = "ON"
$arr[8,0] = 46045.16666666666
$arr[8,1] = 57.06
$arr[8,2] = "forecast"
$arr[8,3] = 46045
$arr[8,4] = "ON"
$arr[9,0] = 46045.1875
$arr[9,1] = 57.09
$arr[9,2] = "forecast"
$arr[9,3] = 46045
$arr[9,4] = "ON"
$arr[10,0] = 46045.20833333334
$arr[10,1] = 57.09
$arr[10,2] = "forecast"
$arr[10,3] = 46045
$arr[10,4] = "ON"
$arr[11,0] = 46045.22916666666
$arr[11,1] = 64.8901
$arr[11,2] = "forecast"
$arr[11,3] = 46045
$arr[11,4] = "ON"
$arr[12,0] = 46045.25
$arr[12,1] = 64.8901
$arr[12,2] = "forecast"
$arr[12,3] = 46045
$arr[12,4] = "ON"
$arr[13,0] = 46045.27083333334
$arr[13,1] = 57.06
$arr[13,2] = "forecast"
$arr[13,3] = 46045
$arr[13,4] = "ON"
$arr[14,0] = 46045.29166666666
$arr[14,1] = 33.46002
$arr[14,2] = "forecast"
$arr[14,3] = 46045
$arr[14,4] = "ON"
$arr[15,0] = 46045.3125
$arr[15,1] = 34.75347
$arr[15,2] = "forecast"
$arr[15,3] = 46045
$arr[15,4] = "ON"
$arr[16,0] = 46045.33333333334
$arr[16,1] = 22.07
$arr[16,2] = "forecast"
$arr[16,3] = 46045
$arr[16,4] = "ON"
$arr[17,0] = 46045.35416666666
$arr[17,1] = 34.16791
$arr[17,2] = "forecast"
$arr[17,3] = 46045
$arr[17,4] = "ON"
$arr[18,0] = 46045.375
$arr[18,1] = 0.51
$arr[18,2] = "forecast"
$arr[18,3] = 46045
$arr[18,4] = "ON"
$arr[19,0] = 46045.39583333334
$arr[19,1] = 0.99404
$arr[19,2] = "forecast"
$arr[19,3] = 46045
$arr[19,4] = "ON"
$arr[20,0] = 46045.41666666666
$arr[20,1] = 36.06
$arr[20,2] = "forecast"
$arr[20,3] = 46045
$arr[20,4] = "ON"
$arr[21,0] = 46045.4375
$arr[21,1] = 36.06
$arr[21,2] = "forecast"
$arr[21,3] = 46045
$arr[21,4] = "ON"
$arr[22,0] = 46045.45833333334
$arr[22,1] = 36.06
$arr[22,2] = "forecast"
$arr[22,3] = 46045
$arr[22,4] = "ON"
$arr[23,0] = 46045.47916666666
$arr[23,1] = 36.06
$arr[23,2] = "forecast"
$arr[23,3] = 46045
$arr[23,4] = "ON"
$arr[24,0] = 46045.5
$arr[24,1] = 36.06
$arr[24,2] = "forecast"
$arr[24,3] = 46045
$arr[24,4] = "ON"
$arr[25,0] = 46045.52083333334
$arr[25,1] = 36.06
$arr[25,2] = "forecast"
$arr[25,3] = 46045
$arr[25,4] = "ON"
$arr[26,0] = 46045.54166666666
$arr[26,1] = 36.06
$arr[26,2] = "forecast"
$arr[26,3] = 46045
$arr[26,4] = "ON"
$arr[27,0] = 46045.5625
$arr[27,1] = 36.06
$arr[27,2] = "forecast"
$arr[27,3] = 46045
$arr[27,4] = "ON"
$arr[28,0] = 46045.58333333334
$arr[28,1] = 36.06
$arr[28,2] = "forecast"
$arr[28,3] = 46045
$arr[28,4] = "ON"
$arr[29,0] = 46045.60416666666
$arr[29,1] = 36.06
$arr[29,2] = "forecast"
$arr[29,3] = 46045
$arr[29,4] = "ON"
$arr[30,0] = 46045.625
$arr[30,1] = 36.06
$arr[30,2] = "forecast"
$arr[30,3] = 46045
$arr[30,4] = "ON"
$arr[31,0] = 46045.64583333334
$arr[31,1] = 22.07
$arr[31,2] = "forecast"
$arr[31,3] = 46045
$arr[31,4] = "ON"
$arr[32,0] = 46045.66666666666
$arr[32,1] = 8.437620000000001
$arr[32,2] = "forecast"
$arr[32,3] = 46045
$arr[32,4] = "OFF"
$arr[33,0] = 46045.6875
$arr[33,1] = 0
$arr[33,2] = "forecast"
$arr[33,3] = 46045
$arr[33,4] = "OFF"
$arr[34,0] = 46045.70833333334
$arr[34,1] = -2.15327
$arr[34,2] = "forecast"
$arr[34,3] = 46045
$arr[34,4] = "OFF"
$arr[35,0] = 46045.72916666666
$arr[35,1] = 4.81867
$arr[35,2] = "forecast"
$arr[35,3] = 46045
$arr[35,4] = "OFF"
$arr[36,0] = 46045.75
$arr[36,1] = 44.27439
$arr[36,2] = "forecast"
$arr[36,3] = 46045
$arr[36,4] = "OFF"
$arr[37,0] = 46045.77083333334
$arr[37,1] = 56.98
$arr[37,2] = "forecast"
$arr[37,3] = 46045
$arr[37,4] = "OFF"
$arr[38,0] = 46045.79166666666
$arr[38,1] = 57.09
$arr[38,2] = "forecast"
$arr[38,3] = 46045
$arr[38,4] = "OFF"
$arr[39,0] = 46045.8125
$arr[39,1] = 59.30547
$arr[39,2] = "forecast"
$arr[39,3] = 46045
$arr[39,4] = "OFF"
$arr[40,0] = 46045.83333333334
$arr[40,1] = 59.37944
$arr[40,2] = "forecast"
$arr[40,3] = 46045
$arr[40,4] = "OFF"
$arr[41,0] = 46045.85416666666
$arr[41,1] = 57.09
$arr[41,2] = "forecast"
$arr[41,3] = 46045
$arr[41,4] = "OFF"
$arr[42,0] = 46045.875
$arr[42,1] = 57.06
$arr[42,2] = "forecast"
$arr[42,3] = 46045
$arr[42,4] = "OFF"
$arr[43,0] = 46045.89583333334
$arr[43,1] = 57.06
$arr[43,2] = "forecast"
$arr[43,3] = 46045
$arr[43,4] = "OFF"
$arr[44,0] = 46045.91666666666
$arr[44,1] = 56.99017
$arr[44,2] = "forecast"
$arr[44,3] = 46045
$arr[44,4] = "OFF"
$arr[45,0] = 46045.9375
$arr[45,1] = 57.06
$arr[45,2] = "forecast"
$arr[45,3] = 46045
$arr[45,4] = "OFF"
$arr[46,0] = 46045.95833333334
$arr[46,1] = 57.06
$arr[46,2] = "forecast"
$arr[46,3] = 46045
$arr[46,4] = "OFF"
$arr[47,0] = 46045.97916666666
$arr[47,1] = 56.98
$arr[47,2] = "forecast"
$arr[47,3] = 46045
$arr[47,4] = "OFF"

$ws2.Range("A2:E49").Value2 = $arr
